# Reorder the block columns and update the corresponding indicator values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order (row 1)
$headers = @("kitchens_1", "kitchens_2", "bedrooms_1", "bedrooms_2", "living_rooms_1", "living_rooms_2")
for ($c = 1; $c -le 6; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# New indicator matrix (rows 2-7), columns aligned to the new headers above
$data = @(
    @(0, 0, 0, 1, 0, 0),
    @(0, 1, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 1),
    @(0, 0, 1, 0, 0, 0),
    @(1, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 1, 0)
)

for ($r = 0; $r -lt 6; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}
